# Timesheet for Week 17 - date corrections
# The "Week of:" date and the seven day-of-week header labels are moved
# forward from the week of 04 May 2014 (Sun 16/02 .. Sat 22/02 - stale
# leftover labels from a prior week that were never updated) to the
# correct week of 27 April 2014 (Sun 27/04 .. Sat 03/05).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Week of:" date field (G8) - serial 41756 = 27 Apr 2014
$ws.Range("G8").Value = 41756

# Day-of-week column headers (A11:A17)
$ws.Range("A11").Value = "Sun 27/04"
$ws.Range("A12").Value = "Mon 28/04"
$ws.Range("A13").Value = "Tue 29/04"
$ws.Range("A14").Value = "Wed 30 /04"
$ws.Range("A15").Value = "Thur 01/05"
$ws.Range("A16").Value = "Fri   02/05"
$ws.Range("A17").Value = "Sat  03/05"
